# Update gh-pages to output generated at 456a3b4
# Applies numeric updates to column F ("想去人数" / interested-attendee count)
# across the four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) -> sheet1.xml
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F10").Value = 796
$wsExhibition.Range("F12").Value = 54
$wsExhibition.Range("F16").Value = 576
$wsExhibition.Range("F20").Value = 835
$wsExhibition.Range("F21").Value = 1156
$wsExhibition.Range("F22").Value = 2844
$wsExhibition.Range("F23").Value = 1376
$wsExhibition.Range("F25").Value = 184
$wsExhibition.Range("F26").Value = 1262
$wsExhibition.Range("F28").Value = 998
$wsExhibition.Range("F30").Value = 2832
$wsExhibition.Range("F31").Value = 565
$wsExhibition.Range("F33").Value = 1379

# Sheet "演出" (Performance) -> sheet2.xml
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F3").Value = 518
$wsPerformance.Range("F13").Value = 3

# Sheet "本地生活" (Local Life) -> sheet3.xml
$wsLocalLife = $wb.Worksheets.Item("本地生活")
$wsLocalLife.Range("F2").Value = 728

# Sheet "全部类型" (All Types) -> sheet4.xml
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 728
$wsAll.Range("F8").Value = 518
$wsAll.Range("F17").Value = 796
$wsAll.Range("F20").Value = 54
$wsAll.Range("F28").Value = 576
$wsAll.Range("F32").Value = 835
$wsAll.Range("F33").Value = 1156
$wsAll.Range("F34").Value = 2844
$wsAll.Range("F35").Value = 1376
$wsAll.Range("F37").Value = 184
$wsAll.Range("F38").Value = 1262
$wsAll.Range("F40").Value = 3
$wsAll.Range("F42").Value = 998
$wsAll.Range("F44").Value = 2832
$wsAll.Range("F45").Value = 565
$wsAll.Range("F47").Value = 1379
